# Set column F ("想去人数") values to 0 for specific rows on sheets
# "展览" (sheet1), "演出" (sheet2), and "全部类型" (sheet4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
for ($r = 2; $r -le 12; $r++) {
    $ws1.Cells.Item($r, 6).Value = 0
}

$ws2 = $wb.Worksheets.Item("演出")
for ($r = 2; $r -le 4; $r++) {
    $ws2.Cells.Item($r, 6).Value = 0
}

$ws4 = $wb.Worksheets.Item("全部类型")
for ($r = 2; $r -le 15; $r++) {
    $ws4.Cells.Item($r, 6).Value = 0
}
